# ======================================================================
# Feat: fotos de lugares adicionadas
#
# For every city in "locais.xlsx" (sheet "Planilha1"), append a new
# "-capa" (cover) image path to the semicolon-separated "fotos" list in
# column C (rows 2-21, one city per row). A handful of cities were still
# pointing at placeholder https://source.unsplash.com/... URLs, and those
# are swapped for real local asset paths (with the new cover appended)
# exactly like the other rows.
#
# Editing C15 (Belem) also drops a leftover "hyperlink-style" underline/
# font that the old placeholder-URL text had been carrying, and Excel
# settles on slightly different auto-fit row heights once the long URL
# text is replaced by the shorter local paths.
# ======================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) "fotos" column (C): append the new cover image for each city ---
$ws.Range("C2").Value = '\assets\images\pontos-tur\rio-de-janeiro.jpg ; \assets\images\pontos-tur\rio-de-janeiro2.jpg ; \assets\images\pontos-tur\rio-de-janeiro3.jpg ; \assets\images\Locais\rio-capa.jpg'
$ws.Range("C3").Value = '\assets\images\Locais\salvador1.jpg ;\assets\images\Locais\salvador2.jpg ; \assets\images\Locais\salvador3.jpg ; \assets\images\Locais\salvador-capa.jpg'
$ws.Range("C4").Value = '\assets\images\Locais\florianopolis1.jpg;\assets\images\Locais\florianopolis2.jpg;\assets\images\Locais\florianopolis3.jpg ; \assets\images\Locais\florianopolis-capa.jpg'
$ws.Range("C5").Value = '\assets\images\Locais\gramado.jpg;\assets\images\Locais\gramado1.webp;\assets\images\Locais\gramado2.jpg ; assets\images\Locais\gramado-capa.jpg'
$ws.Range("C6").Value = '\assets\images\Locais\bonito1.jpg ; \assets\images\Locais\bonito2.JPG ; \assets\images\Locais\bonito3.JPG ; \assets\images\Locais\bonito-capa.webp'
$ws.Range("C7").Value = '\assets\images\Locais\paraty1.jpg ; \assets\images\Locais\paraty2.jpg ; \assets\images\Locais\paraty3.jpg ; \assets\images\Locais\paraty-capa.jpg'
$ws.Range("C8").Value = '\assets\images\Locais\manaus1.jpg ; \assets\images\Locais\manaus2.jpg ; \assets\images\Locais\manaus3.jpg ; \assets\images\Locais\manaus-capa.webp'
$ws.Range("C9").Value = '\assets\images\Locais\foz1.jpg ; \assets\images\Locais\foz2.jpg ; \assets\images\Locais\foz3.jpg ; \assets\images\Locais\foz-capa.jpg'
$ws.Range("C10").Value = '\assets\images\Locais\ouropreto1.jpg ; \assets\images\Locais\ouropreto2.jfif ; \assets\images\Locais\ouropreto3.jpg ; \assets\images\Locais\ouropreto-capa.jpg'
$ws.Range("C11").Value = '\assets\images\Locais\Jericoacoara1.jpg ; \assets\images\Locais\Jericoacoara2.jpg ; \assets\images\Locais\Jericoacoara3.jpg ; \assets\images\Locais\Jericoacoara-capa.jpg'
$ws.Range("C12").Value = '\assets\images\Locais\olinda1.jpg ; \assets\images\Locais\olinda2.jpg ; \assets\images\Locais\olinda3.JPG ; \assets\images\Locais\olinda-capa.jpg'
$ws.Range("C13").Value = '\assets\images\Locais\curitiba1.jpg ; \assets\images\Locais\curitiba2.JPG ; \assets\images\Locais\curitiba3.jpg ; \assets\images\Locais\curitiba-capa.jpg'
$ws.Range("C14").Value = '\assets\images\Locais\natal1.jpg ; \assets\images\Locais\natal2.jpg ; \assets\images\Locais\natal3.jpg ; \assets\images\Locais\natal-capa.webp'
$ws.Range("C15").Value = '\assets\images\Locais\belem1.jpg ; \assets\images\Locais\belem2.jpg ; \assets\images\Locais\belem3.webp ; \assets\images\Locais\belem-capa.webp'
$ws.Range("C16").Value = '\assets\images\Locais\saoluis1.jpg ; \assets\images\Locais\saoluis2.jpg ; \assets\images\Locais\saoluis3.jpg ; \assets\images\Locais\saoluis-capa.jpg'
$ws.Range("C17").Value = '\assets\images\Locais\bh1.jpg ; \assets\images\Locais\bh2.jpg ; \assets\images\Locais\bh3.webp ; \assets\images\Locais\bh-capa.jpg'
$ws.Range("C18").Value = '\assets\images\Locais\recife1.jpg ; \assets\images\Locais\recife3.jfif ; \assets\images\Locais\recife3.webp ; \assets\images\Locais\recife-capa.jpg'
$ws.Range("C19").Value = '\assets\images\Locais\brasilia1.jfif ; \assets\images\Locais\brasilia1.jpg ; \assets\images\Locais\brasilia2.jpg ; \assets\images\Locais\brasilia-capa.jpg'
$ws.Range("C20").Value = '\assets\images\Locais\portoalegre1.png ; \assets\images\Locais\portoalegre2.jpg ; \assets\images\Locais\portoalegre3.jpg ; \assets\images\Locais\portoalegre-capa.jpg'
$ws.Range("C21").Value = '\assets\images\Locais\joaopessoa1.jpg ; \assets\images\Locais\joaopessoa2.jpg ; \assets\images\Locais\joaopessoa3.jpg ; \assets\images\Locais\joaopessoa-capa.jpg'

# --- 2) C15 (Belem) no longer looks like a hyperlink: drop the underline
#        and force the same plain black Calibri font used elsewhere ---
$ws.Range("C15").Font.Underline = $false
$ws.Range("C15").Font.Name = "Calibri"
$ws.Range("C15").Font.Size = 11
$ws.Range("C15").Font.Color = 0
$ws.Range("C15").HorizontalAlignment = -4131

# --- 3) Row heights settle to their new auto-fit values after the edit ---
$ws.Rows.Item(2).RowHeight = 44.25
$ws.Range("A3:A5").RowHeight = 19.5
$ws.Range("A8:A15").RowHeight = 21.75
$ws.Range("A16:A99").RowHeight = 18.75
